$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.398.97'
$ws.Range("E2").Value = '  -1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.048.60'
$ws.Range("E3").Value = '  -2.13%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.95'
$ws.Range("E5").Value = '  -2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  -2.29%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.01'
$ws.Range("E8").Value = '  -4.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.385'
$ws.Range("E9").Value = '  -2.60%  '

$ws.Range("E10").Value = '  +3.46%  '

$ws.Range("E11").Value = '  -2.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.352.67'
$ws.Range("E12").Value = '  -2.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.49'
$ws.Range("E13").Value = '  -5.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.59'
$ws.Range("E14").Value = '  -3.47%  '

$ws.Range("E15").Value = '  -3.56%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.26'
$ws.Range("E16").Value = '  -2.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.052.11'
$ws.Range("E17").Value = '  -1.95%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.302.23'
$ws.Range("E18").Value = '  -1.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.07'
$ws.Range("E19").Value = '  -1.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.82'
$ws.Range("E20").Value = '  -1.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0846'
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.93'
$ws.Range("E22").Value = '  -1.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("E24").Value = '  -0.42%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("E25").Value = '  -5.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.53'
$ws.Range("E26").Value = '  -3.16%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.33'
$ws.Range("E27").Value = '  -1.85%  '

$ws.Range("E28").Value = '  -4.04%  '

$ws.Range("E29").Value = '  -0.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.91'
$ws.Range("E30").Value = '  -3.23%  '

$ws.Range("E31").Value = '  -2.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.54'
$ws.Range("E32").Value = '  -3.89%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.57'
$ws.Range("E33").Value = '  -2.30%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0611'
$ws.Range("E34").Value = '  -3.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.39'
$ws.Range("E35").Value = '  -4.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.82'
$ws.Range("E36").Value = '  +0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.18'
$ws.Range("E38").Value = '  -4.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  +0.76%  '

$ws.Range("E40").Value = '  -6.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.504.68'
$ws.Range("E41").Value = '  +3.58%  '

$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.84'
$ws.Range("E43").Value = '  +0.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.18'
$ws.Range("E44").Value = '  -5.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0931'
$ws.Range("E45").Value = '  -4.50%  '

$ws.Range("E46").Value = '  -3.49%  '

$ws.Range("E47").Value = '  -4.70%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.20'
$ws.Range("E48").Value = '  -1.11%  '

$ws.Range("E49").Value = '  -2.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.238.80'
$ws.Range("E50").Value = '  -2.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.56'
$ws.Range("E51").Value = '  -14.00%  '
